# Commiting scripts (R22 UAT2 - Regression).
#
# The workbook originally has a single worksheet "postRestrictionMarking"
# holding one data row (A2=12105488, ...). This edit:
#   1. Duplicates that worksheet and places the copy BEFORE it, named
#      "Sheet1" - it keeps the original (old) data row and becomes an
#      unselected, "select-all"-ed background tab.
#   2. Overwrites row 2 on the original "postRestrictionMarking" sheet
#      with new values and moves its selection to A2.
#   3. Leaves "postRestrictionMarking" as the active/visible sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Duplicate the existing sheet, insert the copy before it ---------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Copy($ws1)

# NOTE: after Copy(), the COM reference used to invoke it ($ws1) tracks the
# newly-created copy, not the original sheet - so re-resolve both sheets by
# name to get stable handles.
$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "Sheet1"

$orig = $wb.Worksheets.Item("postRestrictionMarking")

# --- 2. Push new row-2 values into the original sheet --------------------
$orig.Range("A2").Value = 17899074
$orig.Range("B2").Value = 42
$orig.Range("C2").Value = 42
$orig.Range("D2").Value = 20230907
$orig.Range("E2").Value = "Blocked by CD - SS Unit"

$orig.Range("A2").Select()

# --- 3. Leftover "select everything" state on the new duplicate tab ------
$newSheet.Activate()
$newSheet.Range("D16").Select()
$newSheet.Cells.Select()

# --- 4. Re-activate the original sheet so it stays the visible tab -------
$orig.Activate()
